$wb = $excel.ActiveWorkbook

$oldGuid = "167eb847-0375-410f-8ab8-e3b9a50a85ec"
$newGuid = "f40ab6d7-f712-4f05-8b15-9649a07f8f04"

# ----------------------------------------------------------------------
# Sheet "Overview"
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# A2: bare file name
$wsOverview.Range("A2").Value = "$newGuid.md"

# B2: "e2e\<guid>.md" -- also carries an external hyperlink whose visible
# text must be refreshed. The link target itself (the .rels Address) is
# unchanged by the diff, only the displayed text changes.
$overviewLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c9148fc816fd1117f0e0f4bc04c2a9ee2c04a83e/e2e/$oldGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewLinkAddress, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md")

# G2: handoff generate date/time
$wsOverview.Range("G2").Value = "2016-10-14 08:11:16"
$wsOverview.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ----------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhCnLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c9148fc816fd1117f0e0f4bc04c2a9ee2c04a83e/e2e/$oldGuid.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhCnLinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md")

# G2: latest handoff xliff file name
$wsZhCn.Range("G2").Value = "$newGuid.34fbb5c45cb61276fec56eb7fb03b16f7c0d0bf3.zh-cn.xlf"

# H2: latest handoff datetime
$wsZhCn.Range("H2").Value = "2016-10-14 08:11:06"
$wsZhCn.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# K2 (Latest Handback DateTime) is untouched by this change, but carries the
# same date-time number format as G2/H2 -- re-assert it so a save of this
# workbook doesn't quietly drop its formatting.
$wsZhCn.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ----------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deDeLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c9148fc816fd1117f0e0f4bc04c2a9ee2c04a83e/e2e/$oldGuid.md"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deDeLinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md")

# G2: latest handoff xliff file name
$wsDeDe.Range("G2").Value = "$newGuid.34fbb5c45cb61276fec56eb7fb03b16f7c0d0bf3.de-de.xlf"

# H2: latest handoff datetime
$wsDeDe.Range("H2").Value = "2016-10-14 08:11:16"
$wsDeDe.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# K2 (Latest Handback DateTime) is untouched by this change, but carries the
# same date-time number format as G2/H2 -- re-assert it so a save of this
# workbook doesn't quietly drop its formatting.
$wsDeDe.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
